# Atualização automática de BOM_JESUS.xlsx
#
# - Remove the "Desarquivamentos Pendentes" sheet (no longer needed).
# - Rename "Paineis DARQ" -> "PAINEIS DARQ".
# - Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO".

$wb = $excel.ActiveWorkbook
[void]($excel.DisplayAlerts = $false)

# Delete the obsolete worksheet entirely.
$wsOld = $wb.Worksheets.Item("Desarquivamentos Pendentes")
[void]$wsOld.Delete()

# Rename the dashboard sheet to all caps.
$wsPaineis = $wb.Worksheets.Item("Paineis DARQ")
$wsPaineis.Name = "PAINEIS DARQ"

# Rename the "Recolhimento x Eliminacao" sheet to all caps (with accent).
$wsRecolhimento = $wb.Worksheets.Item("Recolhimento x Eliminacao")
$wsRecolhimento.Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Keep the dashboard sheet as the active/selected tab.
$wsPaineis.Activate()

[void]($excel.DisplayAlerts = $true)
